$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '42.577.06'
$cell.ClearFormats()
$ws.Range("E2").Value = '  +1.90%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.514.37'
$cell.ClearFormats()
$ws.Range("E3").Value = '  +1.28%  '

$ws.Range("E4").Value = '  -0.33%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '318.03'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +6.31%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '94.49'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +0.54%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("E13").Value = '  -2.26%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '2.900.77'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +1.03%  '

$ws.Range("E15").Value = '  +2.71%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.485.44'
$cell.ClearFormats()
$ws.Range("E16").Value = '  -0.29%  '

$ws.Range("E17").Value = '  +0.05%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '42.649.91'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +1.58%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '12.89'
$cell.ClearFormats()
$ws.Range("E19").Value = '  +2.14%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '6.67'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +5.64%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0956'
$cell.ClearFormats()
$ws.Range("E21").Value = '  -0.01%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '69.29'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -1.80%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '250.45'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +1.78%  '

$ws.Range("E24").Value = '  +2.76%  '

$ws.Range("E25").Value = '  +1.90%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '26.77'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +1.35%  '

$ws.Range("E27").Value = '  +0.21%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.44'
$cell.ClearFormats()
$ws.Range("E28").Value = '  +6.77%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '40.98'
$cell.ClearFormats()
$ws.Range("E29").Value = '  +11.17%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '10.19'
$cell.ClearFormats()
$ws.Range("E30").Value = '  +1.95%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '5.93'
$cell.ClearFormats()
$ws.Range("E31").Value = '  +1.99%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '157.21'
$cell.ClearFormats()
$ws.Range("E32").Value = '  +2.24%  '

$ws.Range("E33").Value = '  +4.48%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '19.06'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +5.50%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.25'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +0.56%  '

$ws.Range("E36").Value = '  +0.77%  '

$ws.Range("E37").Value = '  +0.99%  '

$ws.Range("E38").Value = '  -2.35%  '

$ws.Range("E39").Value = '  +0.06%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '23.37'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -1.82%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '2.29'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +18.57%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("E43").Value = '  +2.52%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '3.31'
$cell.ClearFormats()
$ws.Range("E44").Value = '  -0.16%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.77'
$cell.ClearFormats()
$ws.Range("E45").Value = '  -0.32%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.014.08'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -0.69%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '85.10'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +3.11%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '8.87'
$cell.ClearFormats()
$ws.Range("E48").Value = '  +0.30%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '74.35'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +5.01%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.755.85'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +0.80%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '101.90'
$cell.ClearFormats()
$ws.Range("E51").Value = '  +2.52%  '
